$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a new entry (row 40) to the time-tracking list ---
# Copy the date formatting from the row above (A39) so the new date cell
# reuses the existing "date" cell style instead of creating a new one.
$ws.Range("A39").Copy()
$ws.Range("A40").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = $false

$newDate = Get-Date -Year 2020 -Month 1 -Day 18 -Hour 0 -Minute 0 -Second 0
$ws.Range("A40").Value = $newDate
$ws.Range("B40").Value = "Kursleiter anzeigen; Kursort anzeigen, bearbeiten, anlegen; Teilnehmer anzeigen; offene Posten anzeigen"
$ws.Range("C40").Value = 4

# Recalculate so the SUM/total formulas in column F pick up the new hours.
$excel.CalculateFull()

# --- Update the view to reflect where the user ended up after the edit ---
$aw = $excel.ActiveWindow
$aw.ScrollRow = 31
$aw.ScrollColumn = 1
$ws.Range("C41").Select()
